$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# The "Sum" column (AH) is recomputed for every keyword/technology row so
# that it reflects the updated statistic multiplier: AH = SUM(C:AG) * 11
# (previously the multiplier was 5). This corresponds to the newly added
# statistical analysis across the different technologies used for the plot.
for ($r = 2; $r -le $lastRow; $r++) {
    $sumRange = $ws.Range("C$r`:AG$r")
    $total = 0
    foreach ($cell in $sumRange.Cells) {
        $total += [double]$cell.Value2
    }
    $ws.Range("AH$r").Value2 = $total * 11
}

Write-Host "Recomputed AH2:AH$lastRow"
